$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
# Each Price cell is stored as text in the workbook (e.g. "28.057.50",
# "0.09830"), using dot separators and significant trailing zeros that a
# plain numeric value would not preserve. Force the cell format to Text
# before writing so the literal string is kept exactly as scraped.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.057.50"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.45"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.21"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4965"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3911"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09830"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.109"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.81"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.450"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.55"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.816.72"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.271"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001136"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.23"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.971"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.111.93"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.89"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.71"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.024.25"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.403"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.70"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1054"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.035"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.565"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.609"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06658"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02340"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.902"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2140"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.959"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6200"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.192"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5894"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.696"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.280"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.47"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.938"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.178"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06773"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  -3.71%  "
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("E9").Value = "  +24.63%  "
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("E17").Value = "  +5.16%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  -6.71%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("E51").Value = "  -1.36%  "
